# Weekly update: prepend a new week's Coliflor (Primera/Segunda) price
# observations at the top of the data block (rows 510-511), pushing the
# rest of the historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current first data-continuation
# row (510); Excel shifts everything below down and extends the used
# range/dimension automatically (to R619).
$ws.Rows.Item(510).Insert()
$ws.Rows.Item(510).Insert()

# New row 510: Coliflor, "Primera" quality, week of 2023-11-28.
$ws.Range("A510").Value = 7
$ws.Range("B510").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C510").Value = "Ñuble"
$ws.Range("D510").Value = 45258
$ws.Range("E510").Value = 16
$ws.Range("F510").Value = 100112008
$ws.Range("G510").Value = "Coliflor"
$ws.Range("H510").Value = "Sin especificar"
$ws.Range("I510").Value = "Primera"
$ws.Range("J510").Value = 300
$ws.Range("K510").Value = 1300
$ws.Range("L510").Value = 1300
$ws.Range("M510").Value = 1300
$ws.Range("N510").Value = "`$/unidad"
$ws.Range("O510").Value = "Región del Maule"
$ws.Range("P510").Value = 1300
$ws.Range("Q510").Value = 1
$ws.Range("R510").Value = "Hortaliza"

# New row 511: Coliflor, "Segunda" quality, same week.
$ws.Range("A511").Value = 7
$ws.Range("B511").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C511").Value = "Ñuble"
$ws.Range("D511").Value = 45258
$ws.Range("E511").Value = 16
$ws.Range("F511").Value = 100112008
$ws.Range("G511").Value = "Coliflor"
$ws.Range("H511").Value = "Sin especificar"
$ws.Range("I511").Value = "Segunda"
$ws.Range("J511").Value = 200
$ws.Range("K511").Value = 1000
$ws.Range("L511").Value = 1000
$ws.Range("M511").Value = 1000
$ws.Range("N511").Value = "`$/unidad"
$ws.Range("O511").Value = "Región del Maule"
$ws.Range("P511").Value = 1000
$ws.Range("Q511").Value = 1
$ws.Range("R511").Value = "Hortaliza"
